$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -12.368
$ws.Range("C6").Value = -12.0128
$ws.Range("C7").Value = -13.3205
$ws.Range("E7").Value = 15.6348
$ws.Range("C8").Value = -13.0234
$ws.Range("E11").Value = 16.83919999999999
$ws.Range("E12").Value = 17.84680000000001
$ws.Range("E15").Value = 16.45610000000001
$ws.Range("C16").Value = -14.46969999999999
$ws.Range("C20").Value = -12.16129999999999
$ws.Range("E20").Value = 15.97009999999999
$ws.Range("C21").Value = -13.2772
$ws.Range("E21").Value = 17.0491
$ws.Range("E22").Value = 17.07859999999999
$ws.Range("E23").Value = 16.21039999999999
$ws.Range("C28").Value = -12.2666
$ws.Range("C29").Value = -11.4962
$ws.Range("E29").Value = 17.38000000000001
$ws.Range("C30").Value = -12.44519999999999
$ws.Range("C32").Value = -12.26060000000001
$ws.Range("E34").Value = 17.21330000000001
$ws.Range("C40").Value = -12.428
$ws.Range("E42").Value = 16.50699999999999
$ws.Range("E43").Value = 17.41090000000001
$ws.Range("E44").Value = 16.74399999999999
$ws.Range("E45").Value = 16.5795
$ws.Range("C46").Value = -14.67149999999999
$ws.Range("E46").Value = 16.54970000000001
$ws.Range("E50").Value = 16.57229999999999
$ws.Range("C51").Value = -11.5498
$ws.Range("E51").Value = 17.37270000000001
$ws.Range("C52").Value = -11.289
$ws.Range("C57").Value = -14.27999999999999
$ws.Range("E57").Value = 16.8374
$ws.Range("C59").Value = -13.0107
$ws.Range("C62").Value = -14.304
$ws.Range("E65").Value = 17.41800000000001
$ws.Range("C66").Value = -11.1899
$ws.Range("E66").Value = 17.23780000000001
$ws.Range("E67").Value = 17.11900000000002
$ws.Range("C73").Value = -11.7114
$ws.Range("C74").Value = -11.97230000000001
$ws.Range("C77").Value = -12.35
$ws.Range("E79").Value = 18.50190000000001
$ws.Range("E84").Value = 16.62999999999999
$ws.Range("E87").Value = 16.2983
$ws.Range("C92").Value = -10.9842
$ws.Range("E92").Value = 18.06640000000001
$ws.Range("E97").Value = 16.3821
$ws.Range("C100").Value = -12.6679
